$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the B12:B20 shared-string labels (status descriptions).
$ws.Range("B12").Value = "In elaborazione"
$ws.Range("B13").Value = "Completata con esito positivo"
$ws.Range("B14").Value = "Completata con esito negativo"
$ws.Range("B15").Value = "Annullata"
$ws.Range("B16").Value = "Richieste di estrazione da elaborare"
$ws.Range("B17").Value = "Completata con esito positivo e warning"
$ws.Range("B18").Value = "Notifiche da prendere in carico"
$ws.Range("B19").Value = "Notifiche fruite"
$ws.Range("B20").Value = "Tutte le notifiche "

# Widen column B and mark it explicit (custom width).
$ws.Columns("B").ColumnWidth = 50.6640625

# Reset the B12:B20 rows back to the sheet's default row height / font
# (copy the blank/default formatting from an untouched cell).
$ws.Rows("12:20").RowHeight = 12
$ws.Range("Z1").Copy()
$ws.Range("B12:B20").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Move the selection.
$ws.Range("C1:C1048576").Select()
